$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data for rows 2-93 (ranks 1-92), reflecting the refreshed
# Consolidated_Score scan results (FNO Liquid universe).
$rowsData = @(
    @(1, "TI", 40, 380.9, 585.7, 7, 5.714285714285714, 6, 432.1, "Higher_Probability_Bull_Reversal", "LONG", "31-40"),
    @(2, "DALBHARAT", 36, 2113.55, 2779.36, 7, 5.142857142857143, 6, 2280, "Higher_Probability_Bull_Reversal", "LONG", "31-40"),
    @(3, "UPL", 36, 637.8200000000001, 901.74, 6, 6, 7, 703.8, "Higher_Probability_Bull_Reversal", "LONG", "31-40"),
    @(4, "KCP", 33, 201.92, 263.23, 6, 5.5, 7, 217.25, "Higher_Probability_Bull_Reversal", "LONG", "31-40"),
    @(5, "AMBUJACEM", 31, 580.22, 685.54, 6, 5.166666666666667, 6, 606.55, "Higher_Probability_Bull_Reversal", "LONG", "31-40"),
    @(6, "PAYTM", 31, 952.61, 1200.18, 6, 5.166666666666667, 6, 1014.5, "Higher_Probability_Bull_Reversal", "LONG", "31-40"),
    @(7, "SHK", 30, 222.76, 329.72, 6, 5, 5, 249.5, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(8, "SHYAMMETL", 26, 877, 1074.79, 5, 5.2, 6, 926.45, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(9, "NAUKRI", 25, 1337.31, 1815.27, 4, 6.25, 7, 1456.8, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(10, "SHREECEM", 24, 30348.04, 36215.89, 4, 6, 6, 31815, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(11, "NATIONALUM", 23, 184.73, 226.02, 4, 5.75, 6, 195.05, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(12, "BAJAJ-AUTO", 23, 7974.73, 9513.799999999999, 4, 5.75, 6, 8359.5, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(13, "BORORENEW", 22, 485.56, 924.13, 4, 5.5, 6, 595.2, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(14, "3MINDIA", 22, 29344.46, 34926.61, 4, 5.5, 7, 30740, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(15, "HINDCOPPER", 21, 255.83, 326.72, 4, 5.25, 6, 273.55, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(16, "LATENTVIEW", 21, 402.91, 544.87, 4, 5.25, 6, 438.4, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(17, "GMDCLTD", 21, 367.07, 716.79, 3, 7, 7, 454.5, "Higher_Probability_Bull_Reversal", "LONG", "21-30"),
    @(18, "ICICIBANK", 20, 1399.19, 1643.24, 4, 5, 5, 1460.2, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(19, "GODREJPROP", 20, 2166.42, 3040.75, 4, 5, 5, 2385, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(20, "SWIGGY", 20, 373.49, 546.74, 4, 5, 5, 416.8, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(21, "RHIM", 18, 464.32, 671.25, 3, 6, 6, 516.05, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(22, "AGI", 18, 827.14, 1258.19, 3, 6, 6, 934.9, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(23, "ETERNAL", 17, 252.01, 423.78, 3, 5.666666666666667, 6, 294.95, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(24, "LTF", 17, 196.12, 250.03, 3, 5.666666666666667, 6, 209.6, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(25, "CHEMPLASTS", 17, 405.28, 617.55, 3, 5.666666666666667, 6, 458.35, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(26, "TATACHEM", 17, 908.3200000000001, 1079.25, 3, 5.666666666666667, 7, 951.05, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(27, "NUVOCO", 16, 355.61, 524.5599999999999, 3, 5.333333333333333, 6, 397.85, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(28, "MAXESTATES", 16, 499.53, 692.01, 3, 5.333333333333333, 6, 547.65, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(29, "PRAKASH", 16, 154.72, 268.32, 3, 5.333333333333333, 6, 183.12, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(30, "POWERINDIA", 16, 17948.04, 25415.89, 3, 5.333333333333333, 6, 19815, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(31, "HINDALCO", 15, 656.63, 778.7, 3, 5, 5, 687.15, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(32, "JBCHEPHARM", 15, 1591.1, 1933.9, 3, 5, 5, 1676.8, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(33, "THOMASCOOK", 15, 158.15, 264.36, 3, 5, 5, 184.7, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(34, "CENTENKA", 15, 494.4, 616.79, 3, 5, 5, 525, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(35, "PVRINOX", 15, 944.46, 1303.43, 3, 5, 5, 1034.2, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(36, "PENIND", 15, 230.91, 308.24, 3, 5, 5, 250.24, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(37, "FLUOROCHEM", 15, 3351.99, 4177.63, 3, 5, 5, 3558.4, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(38, "TIMKEN", 15, 3239.83, 4184.5, 3, 5, 5, 3476, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(39, "IIFL", 15, 506.21, 624.96, 3, 5, 5, 535.9, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(40, "ANGELONE", 13, 2574.2, 3540.99, 2, 6.5, 7, 2815.9, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(41, "UBL", 12, 1967.18, 2317.25, 2, 6, 6, 2054.7, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(42, "GODFRYPHLP", 12, 8603.889999999999, 12400.32, 2, 6, 6, 9553, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(43, "SUPREMEIND", 11, 4083.62, 4923.55, 2, 5.5, 6, 4293.6, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(44, "SHRIRAMPPS", 11, 92.58, 122.74, 2, 5.5, 6, 100.12, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(45, "MMTC", 11, 66.08, 87.95, 2, 5.5, 6, 71.55, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(46, "RAYMOND", 11, 677.51, 971.08, 2, 5.5, 6, 750.9, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(47, "PERSISTENT", 11, 5344.59, 7124.23, 2, 6, 6, 5789.5, "Higher_Probability_Bull_Reversal", "LONG", "11-20"),
    @(48, "EXIDEIND", 10, 374.12, 444.84, 2, 5, 5, 391.8, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(49, "MAHLIFE", 10, 358.42, 469.93, 2, 5, 5, 386.3, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(50, "KAJARIACER", 10, 1145.97, 1550.5, 2, 5, 5, 1247.1, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(51, "INDHOTEL", 10, 717.99, 944.8200000000001, 2, 5, 5, 774.7, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(52, "360ONE", 10, 1126.4, 1513.61, 2, 5, 5, 1223.2, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(53, "ASHAPURMIN", 10, 434.15, 1033.36, 2, 5, 5, 583.95, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(54, "TATASTEEL", 10, 155.4, 190.9, 2, 5, 5, 164.27, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(55, "INDIACEM", 10, 325.06, 455.41, 2, 5, 5, 357.65, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(56, "MAPMYINDIA", 10, 1752.79, 2135.64, 2, 5, 5, 1848.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(57, "NLCINDIA", 10, 226.13, 286.68, 2, 5, 5, 241.27, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(58, "LEMONTREE", 10, 145.32, 194.4, 2, 5, 5, 157.59, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(59, "SUNTV", 10, 557.6900000000001, 690.9400000000001, 2, 5, 5, 591, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(60, "M&M", 10, 3050.1, 3834.89, 2, 5, 5, 3246.3, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(61, "VBL", 10, 444.91, 619.27, 2, 5, 5, 488.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(62, "SAMMAANCAP", 10, 121.04, 191.27, 2, 5, 5, 138.6, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(63, "TATAINVEST", 10, 6437.54, 8275.389999999999, 2, 5, 5, 6897, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(64, "SYNGENE", 6, 643.58, 789.26, 1, 6, 6, 680, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(65, "HAVELLS", 6, 1496.21, 1827.37, 1, 6, 6, 1579, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(66, "RBA", 6, 80.98, 102.71, 1, 6, 6, 86.41, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(67, "RAMCOCEM", 6, 1101.98, 1454.05, 1, 6, 6, 1190, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(68, "MASTEK", 5, 2412.24, 3740.07, 1, 5, 5, 2744.2, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(69, "FORTIS", 5, 747.14, 990.77, 1, 5, 5, 808.05, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(70, "DEEPAKFERT", 5, 1497.89, 2032.32, 1, 5, 5, 1631.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(71, "METROBRAND", 5, 1122.79, 1462.43, 1, 5, 5, 1207.7, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(72, "CONCORDBIO", 5, 1750.6, 2512.21, 1, 5, 5, 1941, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(73, "SCHAEFFLER", 5, 4020.86, 5323.01, 1, 5, 5, 4346.4, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(74, "POONAWALLA", 5, 436.59, 536.23, 1, 5, 5, 461.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(75, "ARE&M", 5, 966.17, 1152.7, 1, 5, 5, 1012.8, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(76, "HDFCAMC", 5, 4990.91, 7515.27, 1, 5, 5, 5622, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(77, "SUMICHEM", 5, 544.11, 741.6799999999999, 1, 5, 5, 593.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(78, "PIIND", 5, 3974.05, 4813.84, 1, 5, 5, 4184, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(79, "KRBL", 5, 367.36, 557.92, 1, 5, 5, 415, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(80, "CUMMINSIND", 5, 3427.49, 4159.13, 1, 5, 5, 3610.4, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(81, "CRAFTSMAN", 5, 5801.79, 8712.639999999999, 1, 5, 5, 6529.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(82, "BAJFINANCE", 5, 905.75, 1076.75, 1, 5, 5, 948.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(83, "BANDHANBNK", 5, 170.84, 217.94, 1, 5, 5, 182.61, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(84, "KEI", 5, 3758.98, 4789.05, 1, 5, 5, 4016.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(85, "CGCL", 5, 172.54, 222.42, 1, 5, 5, 185.01, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(86, "BALKRISIND", 5, 2607.32, 3338.05, 1, 5, 5, 2790, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(87, "LLOYDSME", 5, 1443.27, 1852.19, 1, 5, 5, 1545.5, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(88, "BIRLACORPN", 5, 1363.19, 1768.03, 1, 5, 5, 1464.4, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(89, "M&MFIN", 5, 253.36, 305.53, 1, 5, 5, 266.4, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(90, "GREENPLY", 5, 294.41, 500.77, 1, 5, 5, 346, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(91, "JAGSNPHARM", 5, 235.97, 339.3, 1, 5, 5, 261.8, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
    @(92, "JGCHEM", 5, 467.96, 796.11, 1, 5, 5, 550, "Higher_Probability_Bull_Reversal", "LONG", "0-10"),
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $rowsData[$i]
    $sheetRow = $i + 2
    $ws.Cells.Item($sheetRow, 1).Value = $r[0]
    $ws.Cells.Item($sheetRow, 2).Value = $r[1]
    $ws.Cells.Item($sheetRow, 3).Value = $r[2]
    $ws.Cells.Item($sheetRow, 4).Value = $r[3]
    $ws.Cells.Item($sheetRow, 5).Value = $r[4]
    $ws.Cells.Item($sheetRow, 6).Value = $r[5]
    $ws.Cells.Item($sheetRow, 7).Value = $r[6]
    $ws.Cells.Item($sheetRow, 8).Value = $r[7]
    $ws.Cells.Item($sheetRow, 9).Value = $r[8]
    $ws.Cells.Item($sheetRow, 10).Value = $r[9]
    $ws.Cells.Item($sheetRow, 11).Value = $r[10]
    $ws.Cells.Item($sheetRow, 12).Value = $r[11]
}

# The refreshed scan has one fewer row than before (92 vs 93 data rows),
# so remove the now-unused last row (previously row 94 / ELGIEQUIP).
$ws.Rows.Item(94).Delete()

